$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the "Gracia María Escalante Iraheta 20255999" paragraph
#    (one of the student-name paragraphs at the top of the document).
# ------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text.Contains("Escalante Iraheta 20255999")) {
        $para.Range.Delete()
        break
    }
}

# ------------------------------------------------------------------
# 2) Turn " (¡BANANA!)" into " ¡BANANA!" -- drop the parentheses that
#    wrap the exclamation, keep the "¡BANANA!" text itself untouched.
# ------------------------------------------------------------------
$rng = $d.Content
$rng.Start = 0
$found = $rng.Find.Execute("(¡BANANA!)")
if ($found) {
    # $rng now spans exactly "(¡BANANA!)" -- shrink it to the leading
    # "(" and trailing ")" characters and delete only those, so the
    # "¡BANANA!" text itself is left completely intact.
    $openParen = $d.Range($rng.Start, $rng.Start + 1)
    $closeParen = $d.Range($rng.End - 1, $rng.End)
    $closeParen.Delete()
    $openParen.Delete()
}

# ------------------------------------------------------------------
# 3) Insert a new centered, non-bold paragraph reading
#    "Documentación del código" right after the paragraph that ends
#    in "...¡BANANA!".
# ------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text.Contains("BANANA!")) {
        $para.Range.InsertParagraphAfter()
        $newPara = $d.Paragraphs.Item($i + 1)
        $newPara.Range.InsertBefore("Documentación del código")
        $newPara.Alignment = 1
        break
    }
}
